# tratando dados faltantes na coluna estado
# For rows where "estado" (column K) is "Não informado" and the job is
# fully remote ("modalidade" column L == "Remoto"), we now know the
# vacancy is open to the whole country, so we set estado to "Todo o Brasil".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $estado = $ws.Cells.Item($r, 11).Value2      # column K = estado
    $modalidade = $ws.Cells.Item($r, 12).Value2  # column L = modalidade

    if ($estado -eq "Não informado" -and $modalidade -eq "Remoto") {
        $ws.Cells.Item($r, 11).Value = "Todo o Brasil"
    }
}
